{"js": "// Append two new sub-bullet (\"ilvl 1\") items to the end of the \"High\"\n// priority list, right after the existing\n// \"Implement \"remember me\" and \"reset password\" features\" item.\n// Both new paragraphs must inherit that item's list/paragraph formatting\n// (ListParagraph style, ilvl 1 / numId 3 bullet, justified) and its\n// \"HTML Code\" run character style, so we simply insert new paragraphs\n// right after the last paragraph in the document body \u2014 Word carries the\n// source paragraph's formatting onto the freshly inserted paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The target paragraph (\"Implement \u201cremember me\u201d and \u201creset password\u201d\n// features\") is the very last paragraph in the document.\nlet lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newItems = [\n  \"User should be redirected to the desired page after successful login when initially access to that page was denied\",\n  \"Do we need to allow to delete company\\u2019s last admin, otherwise the company cannot be deleted\"\n];\n\nfor (const text of newItems) {\n  lastParagraph = lastParagraph.insertParagraph(text, \"After\");\n  await context.sync();\n}\n", "ps1": "# Append two new sub-bullet (\"ilvl 1\") items to the end of the \"High\"\n# priority list, right after the existing\n# \"Implement \"remember me\" and \"reset password\" features\" item.\n# Both new paragraphs must inherit that item's list/paragraph formatting\n# (ListParagraph style, ilvl 1 / numId 3 bullet, justified) and its\n# \"HTML Code\" run character style, so we simply insert new paragraphs\n# right after the last paragraph in the document \u2014 Word carries the\n# source paragraph's formatting onto the freshly inserted paragraph.\n\n$d = $word.ActiveDocument\n\n$rsquo = [char]0x2019\n\n$newItems = @(\n    \"User should be redirected to the desired page after successful login when initially access to that page was denied\",\n    \"Do we need to allow to delete company${rsquo}s last admin, otherwise the company cannot be deleted\"\n)\n\n$lastParagraph = $d.Paragraphs.Last\n\nforeach ($text in $newItems) {\n    $lastParagraph.Range.InsertParagraphAfter()\n    $lastParagraph = $d.Paragraphs.Last\n    $lastParagraph.Range.Text = $text\n}\n"}
